# Generate Report for Handoff
#
# Refreshes the "Latest Handoff Date(time)" values for the files that were
# just (re-)handed off, on each of the three report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D is "Latest Handoff Date" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D14").Value = "2016-11-18 04:11:14"
$overview.Range("D15").Value = "2016-11-18 04:11:14"

# --- zh-cn sheet: column E is "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $zhcnRows) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-18 04:11:11"
}

# --- de-de sheet: column E is "Latest Handoff Datetime" ---
$dede = $wb.Worksheets.Item("de-de")
$dedeRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $dedeRows) {
    $dede.Cells.Item($r, 5).Value = "2016-03-18 04:11:14"
}
